# Add a new "2022" column (S) to the report, mirroring the formatting of
# the existing "2021" column (R), then set the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / font / alignment / borders) from the
# 2021 column (R4:R14) into the new 2022 column (S4:S14) so the new cells
# look consistent with the rest of the table.
$ws.Range("R4:R14").Copy() | Out-Null
$ws.Range("S4:S14").PasteSpecial(-4122) | Out-Null

# New header year.
$ws.Range("S4").Value = 2022

# New data values (all zero, no observation yet for 2022).
$ws.Range("S5").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("S14").Value = 0

# Match the author's final selection (saved in the sheet view).
$ws.Range("R17").Select()
